$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") "63.407.68"
Set-TextCell $ws.Range("E2") "  +1.81%  "

# Row 3
Set-TextCell $ws.Range("D3") "3.173.02"
Set-TextCell $ws.Range("E3") "  -0.48%  "

# Row 4
Set-TextCell $ws.Range("E4") "  +0.05%  "

# Row 5
Set-TextCell $ws.Range("D5") "601.70"
Set-TextCell $ws.Range("E5") "  +2.60%  "

# Row 6
Set-TextCell $ws.Range("D6") "136.09"
Set-TextCell $ws.Range("E6") "  +0.59%  "

# Row 7
Set-TextCell $ws.Range("E7") "  +0.05%  "

# Row 8
Set-TextCell $ws.Range("D8") "3.170.78"
Set-TextCell $ws.Range("E8") "  -0.51%  "

# Row 9
Set-TextCell $ws.Range("E9") "  +1.99%  "

# Row 10
Set-TextCell $ws.Range("E10") "  +0.83%  "

# Row 11
Set-TextCell $ws.Range("E11") "  +2.61%  "

# Row 12
Set-TextCell $ws.Range("D12") "0.455"
Set-TextCell $ws.Range("E12") "  +0.66%  "

# Row 13
Set-TextCell $ws.Range("E13") "  +2.17%  "

# Row 14
Set-TextCell $ws.Range("D14") "34.97"
Set-TextCell $ws.Range("E14") "  +5.23%  "

# Row 15
Set-TextCell $ws.Range("D15") "3.695.69"
Set-TextCell $ws.Range("E15") "  -0.27%  "

# Row 16
Set-TextCell $ws.Range("D16") "0.121"
Set-TextCell $ws.Range("E16") "  +1.27%  "

# Row 17
Set-TextCell $ws.Range("D17") "3.170.85"
Set-TextCell $ws.Range("E17") "  -0.23%  "

# Row 18
Set-TextCell $ws.Range("D18") "63.374.78"
Set-TextCell $ws.Range("E18") "  +1.65%  "

# Row 19
Set-TextCell $ws.Range("D19") "6.60"
Set-TextCell $ws.Range("E19") "  +0.07%  "

# Row 20
Set-TextCell $ws.Range("D20") "461.77"
Set-TextCell $ws.Range("E20") "  +1.20%  "

# Row 21
Set-TextCell $ws.Range("E21") "  -0.07%  "

# Row 22
Set-TextCell $ws.Range("D22") "0.698"
Set-TextCell $ws.Range("E22") "  -1.10%  "

# Row 23
Set-TextCell $ws.Range("E23") "  +0.55%  "

# Row 24
Set-TextCell $ws.Range("D24") "13.35"
Set-TextCell $ws.Range("E24") "  -0.51%  "

# Row 25
Set-TextCell $ws.Range("D25") "83.27"
Set-TextCell $ws.Range("E25") "  +0.78%  "

# Row 27
Set-TextCell $ws.Range("E27") "  +0.75%  "

# Row 28
Set-TextCell $ws.Range("D28") "0.999"
Set-TextCell $ws.Range("E28") "  +0.05%  "

# Row 29
Set-TextCell $ws.Range("E29") "  +3.66%  "

# Row 30
Set-TextCell $ws.Range("B30") "NEARProtocol"
Set-TextCell $ws.Range("C30") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell $ws.Range("D30") "6.80"
Set-TextCell $ws.Range("E30") "  -1.54%  "

# Row 31
Set-TextCell $ws.Range("B31") "RenderToken"
Set-TextCell $ws.Range("C31") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws.Range("D31") "7.71"
Set-TextCell $ws.Range("E31") "  -1.65%  "

# Row 32
Set-TextCell $ws.Range("D32") "27.15"
Set-TextCell $ws.Range("E32") "  -0.53%  "

# Row 33
Set-TextCell $ws.Range("E33") "  -1.93%  "

# Row 34
Set-TextCell $ws.Range("E34") "  +1.75%  "

# Row 35
Set-TextCell $ws.Range("E35") "  -1.91%  "

# Row 36
Set-TextCell $ws.Range("D36") "5.92"
Set-TextCell $ws.Range("E36") "  +1.96%  "

# Row 37
Set-TextCell $ws.Range("D37") "0.0₃0735"
Set-TextCell $ws.Range("E37") "  +6.36%  "

# Row 38
Set-TextCell $ws.Range("D38") "51.24"
Set-TextCell $ws.Range("E38") "  +0.36%  "

# Row 39
Set-TextCell $ws.Range("D39") "0.0392"
Set-TextCell $ws.Range("E39") "  +1.74%  "

# Row 40
Set-TextCell $ws.Range("D40") "8.15"
Set-TextCell $ws.Range("E40") "  +1.41%  "

# Row 41
Set-TextCell $ws.Range("B41") "dogwifhat"
Set-TextCell $ws.Range("C41") "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell $ws.Range("D41") "2.64"
Set-TextCell $ws.Range("E41") "  +0.26%  "

# Row 42
Set-TextCell $ws.Range("B42") "Kaspa"
Set-TextCell $ws.Range("C42") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell $ws.Range("D42") "0.112"
Set-TextCell $ws.Range("E42") "  -0.61%  "

# Row 43
Set-TextCell $ws.Range("D43") "392.22"
Set-TextCell $ws.Range("E43") "  -4.72%  "

# Row 44
Set-TextCell $ws.Range("D44") "2.802.47"
Set-TextCell $ws.Range("E44") "  -5.04%  "

# Row 45
Set-TextCell $ws.Range("E45") "  +0.28%  "

# Row 46
Set-TextCell $ws.Range("D46") "36.11"
Set-TextCell $ws.Range("E46") "  -0.02%  "

# Row 47
Set-TextCell $ws.Range("D47") "2.13"
Set-TextCell $ws.Range("E47") "  -1.37%  "

# Row 49
Set-TextCell $ws.Range("E49") "  +2.52%  "

# Row 50
Set-TextCell $ws.Range("D50") "25.20"
Set-TextCell $ws.Range("E50") "  -1.32%  "

# Row 51
Set-TextCell $ws.Range("E51") "  +0.83%  "

